$d = $word.ActiveDocument

# 1. Remove the childless-tag paragraph "PUMP:HRD:0000" entirely (whole
#    paragraph, including its paragraph mark).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*PUMP:HRD:0000*") {
        $p.Range.Delete()
    }
}

# 2. Insert a new paragraph "PUMP:HTR:200" right after "PUMP:HTR:1500".
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*PUMP:HTR:1500*") {
        $targetIndex = $i
        break
    }
}
$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "PUMP:HTR:200"

# 3. Remove the old "PUMPHTR:200" paragraph (the tag is now represented
#    by the newly inserted "PUMP:HTR:200" paragraph above).
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*PUMPHTR:200*") {
        $p.Range.Delete()
    }
}
